# Edit script: add cardiff_properties and barry_properties sheets,
# and clear the stray empty B9/C9 cells on the newcastle_upon_tyne_properties sheet.

$wb = $excel.ActiveWorkbook
$wsNewcastle = $wb.Worksheets.Item(1)

# --- Clear the two stray empty cells on row 9 (Quayside Stay...) ---
$wsNewcastle.Range("B9").ClearContents()
$wsNewcastle.Range("C9").ClearContents()

# Keep a handle on the header formatting used on the first sheet so the
# new sheets' header rows can reuse the same bold/bordered style.
$headerFormat = $wsNewcastle.Range("A1:D1")

# --- Add cardiff_properties sheet after the last existing sheet ---
$wsCardiff = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsCardiff.Name = "cardiff_properties"
$wsCardiff.PageSetup.LeftMargin = 54
$wsCardiff.PageSetup.RightMargin = 54
$wsCardiff.PageSetup.TopMargin = 72
$wsCardiff.PageSetup.BottomMargin = 72
$wsCardiff.PageSetup.HeaderMargin = 36
$wsCardiff.PageSetup.FooterMargin = 36

$wsCardiff.Cells.Item(1, 1).Value = "name"
$wsCardiff.Cells.Item(1, 2).Value = "average_rating"
$wsCardiff.Cells.Item(1, 3).Value = "number_of_ratings"
$wsCardiff.Cells.Item(1, 4).Value = "price_per_night"

$wsCardiff.Cells.Item(2, 1).Value = 'A lovely private...'
$wsCardiff.Cells.Item(2, 2).Value = 4.86
$wsCardiff.Cells.Item(2, 3).Value = 7
$wsCardiff.Cells.Item(2, 4).Value = 53
$wsCardiff.Cells.Item(3, 1).Value = 'Room in Cardiff ...'
$wsCardiff.Cells.Item(3, 4).Value = 47
$wsCardiff.Cells.Item(4, 1).Value = 'City Centre Retr...'
$wsCardiff.Cells.Item(4, 2).Value = 4.87
$wsCardiff.Cells.Item(4, 3).Value = 366
$wsCardiff.Cells.Item(4, 4).Value = 123
$wsCardiff.Cells.Item(5, 1).Value = 'Warm & welcoming...'
$wsCardiff.Cells.Item(5, 2).Value = 4.97
$wsCardiff.Cells.Item(5, 3).Value = 951
$wsCardiff.Cells.Item(5, 4).Value = 76
$wsCardiff.Cells.Item(6, 1).Value = 'Cosy & Central S...'
$wsCardiff.Cells.Item(6, 2).Value = 4.79
$wsCardiff.Cells.Item(6, 3).Value = 53
$wsCardiff.Cells.Item(6, 4).Value = 70
$wsCardiff.Cells.Item(7, 1).Value = 'Central & Modern...'
$wsCardiff.Cells.Item(7, 2).Value = 4.94
$wsCardiff.Cells.Item(7, 3).Value = 84
$wsCardiff.Cells.Item(7, 4).Value = 59
$wsCardiff.Cells.Item(8, 1).Value = 'Detached, indepe...'
$wsCardiff.Cells.Item(8, 2).Value = 4.97
$wsCardiff.Cells.Item(8, 3).Value = 380
$wsCardiff.Cells.Item(8, 4).Value = 73
$wsCardiff.Cells.Item(9, 1).Value = '5 mins to Centre...'
$wsCardiff.Cells.Item(9, 2).Value = 4.79
$wsCardiff.Cells.Item(9, 3).Value = 192
$wsCardiff.Cells.Item(9, 4).Value = 90
$wsCardiff.Cells.Item(10, 1).Value = 'Single room in g...'
$wsCardiff.Cells.Item(10, 2).Value = 4.77
$wsCardiff.Cells.Item(10, 3).Value = 329
$wsCardiff.Cells.Item(10, 4).Value = 26
$wsCardiff.Cells.Item(11, 1).Value = 'Central & Modern...'
$wsCardiff.Cells.Item(11, 2).Value = 4.89
$wsCardiff.Cells.Item(11, 3).Value = 94
$wsCardiff.Cells.Item(11, 4).Value = 55
$wsCardiff.Cells.Item(12, 1).Value = 'The Little Lake ...'
$wsCardiff.Cells.Item(12, 2).Value = 4.88
$wsCardiff.Cells.Item(12, 3).Value = 739
$wsCardiff.Cells.Item(12, 4).Value = 74
$wsCardiff.Cells.Item(13, 1).Value = 'Spacious Detache...'
$wsCardiff.Cells.Item(13, 2).Value = 4.99
$wsCardiff.Cells.Item(13, 3).Value = 202
$wsCardiff.Cells.Item(13, 4).Value = 104
$wsCardiff.Cells.Item(14, 1).Value = 'Cosy Victorian h...'
$wsCardiff.Cells.Item(14, 2).Value = 4.9
$wsCardiff.Cells.Item(14, 3).Value = 236
$wsCardiff.Cells.Item(14, 4).Value = 44
$wsCardiff.Cells.Item(15, 1).Value = 'Cosy cabin style...'
$wsCardiff.Cells.Item(15, 2).Value = 4.88
$wsCardiff.Cells.Item(15, 3).Value = 505
$wsCardiff.Cells.Item(15, 4).Value = 38
$wsCardiff.Cells.Item(16, 1).Value = 'Double Room / Ow...'
$wsCardiff.Cells.Item(16, 2).Value = 4.99
$wsCardiff.Cells.Item(16, 3).Value = 543
$wsCardiff.Cells.Item(16, 4).Value = 81
$wsCardiff.Cells.Item(17, 1).Value = 'Compact Tiny Taf...'
$wsCardiff.Cells.Item(17, 2).Value = 4.97
$wsCardiff.Cells.Item(17, 3).Value = 170
$wsCardiff.Cells.Item(17, 4).Value = 95
$wsCardiff.Cells.Item(18, 1).Value = '(B) Private En-s...'
$wsCardiff.Cells.Item(18, 2).Value = 4.91
$wsCardiff.Cells.Item(18, 3).Value = 284
$wsCardiff.Cells.Item(18, 4).Value = 58

$headerFormat.Copy()
$wsCardiff.Range("A1:D1").PasteSpecial(-4122)

# --- Add barry_properties sheet after cardiff_properties ---
$wsBarry = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsBarry.Name = "barry_properties"
$wsBarry.PageSetup.LeftMargin = 54
$wsBarry.PageSetup.RightMargin = 54
$wsBarry.PageSetup.TopMargin = 72
$wsBarry.PageSetup.BottomMargin = 72
$wsBarry.PageSetup.HeaderMargin = 36
$wsBarry.PageSetup.FooterMargin = 36

$wsBarry.Cells.Item(1, 1).Value = "name"
$wsBarry.Cells.Item(1, 2).Value = "average_rating"
$wsBarry.Cells.Item(1, 3).Value = "number_of_ratings"
$wsBarry.Cells.Item(1, 4).Value = "price_per_night"

$wsBarry.Cells.Item(2, 1).Value = 'The Annex @ Broo...'
$wsBarry.Cells.Item(2, 2).Value = 4.95
$wsBarry.Cells.Item(2, 3).Value = 153
$wsBarry.Cells.Item(2, 4).Value = 53
$wsBarry.Cells.Item(3, 1).Value = '"Y Sied"-quiet, ...'
$wsBarry.Cells.Item(3, 2).Value = 4.98
$wsBarry.Cells.Item(3, 3).Value = 206
$wsBarry.Cells.Item(3, 4).Value = 111
$wsBarry.Cells.Item(4, 1).Value = 'Double bed w/ en...'
$wsBarry.Cells.Item(4, 2).Value = 4.98
$wsBarry.Cells.Item(4, 3).Value = 64
$wsBarry.Cells.Item(4, 4).Value = 66
$wsBarry.Cells.Item(5, 1).Value = 'Comfy small room...'
$wsBarry.Cells.Item(5, 2).Value = 4.71
$wsBarry.Cells.Item(5, 3).Value = 340
$wsBarry.Cells.Item(5, 4).Value = 37
$wsBarry.Cells.Item(6, 1).Value = 'Crow''s Nest Barr...'
$wsBarry.Cells.Item(6, 2).Value = 4.94
$wsBarry.Cells.Item(6, 3).Value = 51
$wsBarry.Cells.Item(6, 4).Value = 102
$wsBarry.Cells.Item(7, 1).Value = 'Cosy Gladstone...'
$wsBarry.Cells.Item(7, 2).Value = 4.72
$wsBarry.Cells.Item(7, 3).Value = 250
$wsBarry.Cells.Item(7, 4).Value = 40
$wsBarry.Cells.Item(8, 1).Value = 'Sea view, entire...'
$wsBarry.Cells.Item(8, 2).Value = 4.77
$wsBarry.Cells.Item(8, 3).Value = 64
$wsBarry.Cells.Item(8, 4).Value = 138
$wsBarry.Cells.Item(9, 1).Value = 'Lovely light sum...'
$wsBarry.Cells.Item(9, 2).Value = 4.98
$wsBarry.Cells.Item(9, 3).Value = 145
$wsBarry.Cells.Item(9, 4).Value = 102
$wsBarry.Cells.Item(10, 1).Value = 'Bright Seaside H...'
$wsBarry.Cells.Item(10, 2).Value = 4.96
$wsBarry.Cells.Item(10, 3).Value = 212
$wsBarry.Cells.Item(10, 4).Value = 181
$wsBarry.Cells.Item(11, 1).Value = 'Vale View (Barry...'
$wsBarry.Cells.Item(11, 2).Value = 4.97
$wsBarry.Cells.Item(11, 3).Value = 29
$wsBarry.Cells.Item(11, 4).Value = 164
$wsBarry.Cells.Item(12, 1).Value = 'Detached house o...'
$wsBarry.Cells.Item(12, 2).Value = 5
$wsBarry.Cells.Item(12, 3).Value = 6
$wsBarry.Cells.Item(12, 4).Value = 48
$wsBarry.Cells.Item(13, 1).Value = 'The 19th Dock...'
$wsBarry.Cells.Item(13, 2).Value = 5
$wsBarry.Cells.Item(13, 3).Value = 9
$wsBarry.Cells.Item(13, 4).Value = 145
$wsBarry.Cells.Item(14, 1).Value = 'luxury 2 bedroom...'
$wsBarry.Cells.Item(14, 2).Value = 4.99
$wsBarry.Cells.Item(14, 3).Value = 83
$wsBarry.Cells.Item(14, 4).Value = 178
$wsBarry.Cells.Item(15, 1).Value = 'The Little Blue ...'
$wsBarry.Cells.Item(15, 2).Value = 4.88
$wsBarry.Cells.Item(15, 3).Value = 8
$wsBarry.Cells.Item(15, 4).Value = 199
$wsBarry.Cells.Item(16, 1).Value = 'Large detached h...'
$wsBarry.Cells.Item(16, 2).Value = 4.9
$wsBarry.Cells.Item(16, 3).Value = 49
$wsBarry.Cells.Item(16, 4).Value = 60
$wsBarry.Cells.Item(17, 1).Value = 'Room - Garden Vi...'
$wsBarry.Cells.Item(17, 2).Value = 5
$wsBarry.Cells.Item(17, 3).Value = 3
$wsBarry.Cells.Item(17, 4).Value = 88
$wsBarry.Cells.Item(18, 1).Value = 'Ty Hapus is a 3-...'
$wsBarry.Cells.Item(18, 2).Value = 4.97
$wsBarry.Cells.Item(18, 3).Value = 193
$wsBarry.Cells.Item(18, 4).Value = 148

$headerFormat.Copy()
$wsBarry.Range("A1:D1").PasteSpecial(-4122)

$wsNewcastle.Select()
